# Data Scrape and Validation Project Plan.xlsx -- apply "Add files via upload" edit
# Updates the "Project Plan" sheet: fills in completion dates/percentages for a
# second pass of the Legislative Activity Data checklist plus the Deliverable and
# Project Closeout rows, relabels a couple of checklist items, adds "Available in
# GitHub" notes, and moves the active selection back up to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Plan")

# --- Deliverable rows 26/27: add a Notes entry -----------------------------
$ws.Range("D26").Value = "Available in GitHub"
$ws.Range("D27").Value = "Available in GitHub"

# --- Legislative Activity Data checklist (rows 60-64), first pass ----------
$ws.Range("A60").Value = "Legislative Activity Data"
$ws.Range("B63").Value = 45376
$ws.Range("C63").Value = 1
$ws.Range("A64").Value = "Update dataset as needed"
$ws.Range("B64").Value = 45376
$ws.Range("C64").Value = 1

# --- Legislative Activity Data checklist (rows 65-69), second pass ---------
$ws.Range("A65").Value = "Legislative Activity Data"
$ws.Range("B66").Value = 45376
$ws.Range("C66").Value = 1
$ws.Range("A67").Value = "Verify data integrity and check for outliers"
$ws.Range("B67").Value = 45376
$ws.Range("C67").Value = 1
$ws.Range("A68").Value = "Verify totals match detail lines"
$ws.Range("B68").Value = 45378
$ws.Range("C68").Value = 1
$ws.Range("A69").Value = "Update dataset as needed"
$ws.Range("B69").Value = 45378
$ws.Range("C69").Value = 1

# --- Deliverable rows 70-72 -------------------------------------------------
$ws.Range("A70").Value = "Deliverable: Final Legislative Activity dataset in Excel format"
$ws.Range("B70").Value = 45378
$ws.Range("C70").Value = 1
$ws.Range("D70").Value = "Available in GitHub"

$ws.Range("A71").Value = "Deliverable: Final Confirmation dataset in Excel format"
$ws.Range("B71").Value = 45378
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = "Available in GitHub"

$ws.Range("A72").Value = "Deliverable: Document data issues and updates in final presentation"
$ws.Range("B72").Value = 45378
$ws.Range("C72").Value = 1
$ws.Range("C72").NumberFormat = "0%"
$ws.Range("D72").Value = "Available in GitHub"

# --- Project Closeout rows 75-76 -------------------------------------------
$ws.Range("B75").Value = 45380
$ws.Range("C75").Value = 1
$ws.Range("C75").NumberFormat = "0%"
$ws.Range("D75").Value = "Available in GitHub"

$ws.Range("B76").Value = 45380
$ws.Range("C76").Value = 1
$ws.Range("C76").NumberFormat = "0%"
$ws.Range("D76").Value = "Available in GitHub"

# --- Restore view state: scroll back to top, select A5 ---------------------
$ws.Range("A5").Select()
